$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append the new daily allocation row (10/26/2025) right after the last
# existing data row (row 54), extending the sheet to row 55.
$row = 55

# Column A: date kept as literal text (matching the existing rows, which
# store dates as plain strings rather than Excel date serials). Apply a
# Text number format before writing so Excel doesn't auto-convert the
# "10/26/2025" string into a date value, then clear the formatting so the
# new cell doesn't end up with a lingering explicit style (matching the
# unstyled data cells above it).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "10/26/2025"
$ws.Cells.Item($row, 1).ClearFormats()

# Columns B and C: the BTC / KAS allocation fractions for that date.
$ws.Cells.Item($row, 2).Value = 0.1833298628465994
$ws.Cells.Item($row, 3).Value = 0.8166701371534006
